# Update cryptos list values (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "27.631.41"
$ws.Range("E2").Value = "  -1.29%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.632.60"
$ws.Range("E3").Value = "  -0.61%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.10%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.38"
$ws.Range("E5").Value = "  -0.68%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  -0.77%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.12%  "

# Row 8 - Solana
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.08"
$ws.Range("E8").Value = "  -0.56%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +0.02%  "

# Row 10 - Dogecoin
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0611"
$ws.Range("E10").Value = "  -0.24%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  -3.30%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.864.03"
$ws.Range("E12").Value = "  -0.59%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.630.40"
$ws.Range("E13").Value = "  +0.04%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  -0.32%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  -0.16%  "

# Row 16 - Litecoin
$ws.Range("E16").Value = "  +0.64%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "27.610.09"
$ws.Range("E17").Value = "  -1.27%  "

# Row 18 - BitcoinCash
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "229.39"
$ws.Range("E18").Value = "  -1.54%  "

# Row 19 - ShibaInu
$ws.Range("D19").Value = "0.0₃0719"
$ws.Range("E19").Value = "  -0.63%  "

# Row 20 - Chainlink
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.55"
$ws.Range("E20").Value = "  -1.37%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  -0.16%  "

# Row 22 - Avalanche
$ws.Range("E22").Value = "  +6.39%  "

# Row 23 - Uniswap
$ws.Range("E23").Value = "  +1.35%  "

# Row 24 - Toncoin
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.14"
$ws.Range("E24").Value = "  +2.76%  "

# Row 25 - Monero
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "149.15"
$ws.Range("E25").Value = "  -0.59%  "

# Row 26 - Cosmos
$ws.Range("E26").Value = "  -1.13%  "

# Row 27 - Stellar
$ws.Range("E27").Value = "  -0.84%  "

# Row 28 - EthereumClassic
$ws.Range("E28").Value = "  -0.32%  "

# Row 29 - BinanceUSD
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.12%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  -0.82%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  -1.28%  "

# Row 33 - Maker
$ws.Range("D33").Value = "1.462.20"
$ws.Range("E33").Value = "  -0.74%  "

# Row 34 - InternetComputer(DFINITY)
$ws.Range("E34").Value = "  -0.57%  "

# Row 35 - LidoDAOToken
$ws.Range("E35").Value = "  -0.31%  "

# Row 36 - HuobiToken
$ws.Range("E36").Value = "  -1.81%  "

# Row 37 - ARBITRUM
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.876"
$ws.Range("E37").Value = "  -0.47%  "

# Row 38 <-> Row 40 swap (VeChain and TrustWalletToken swap ranking order)
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0167"
$ws.Range("E38").Value = "  -0.44%  "

# Row 39 - ImmutableX
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.557"
$ws.Range("E39").Value = "  -1.92%  "

$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.920"
$ws.Range("E40").Value = "  -0.31%  "

# Row 41 - Aave
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "68.95"
$ws.Range("E41").Value = "  -0.71%  "

# Row 42 - PaxDollar
$ws.Range("E42").Value = "  -0.07%  "

# Row 43 - WEMIXToken
$ws.Range("E43").Value = "  +0.44%  "

# Row 44 - mCoin
$ws.Range("E44").Value = "  -0.12%  "

# Row 45 <-> Row 46 swap (FraxShare and MXToken swap ranking order)
$ws.Range("B45").Value = "MXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.21"
$ws.Range("E45").Value = "  -1.12%  "

$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.38"
$ws.Range("E46").Value = "  -0.39%  "

# Row 47 - RocketPoolETH
$ws.Range("D47").Value = "1.773.70"

# Row 48 - RenderToken
$ws.Range("E48").Value = "  +2.47%  "

# Row 49 - Quant
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "87.51"
$ws.Range("E49").Value = "  +1.62%  "

# Row 50 - BabyDogeCoin
$ws.Range("E50").Value = "  -0.87%  "

# Row 51 - Algorand
$ws.Range("E51").Value = "  +0.25%  "
